# Update results: shift trial index (column B) by +1 for every trial row
# (rows 2-41) across all participant worksheets, and fix a few prediction
# (column C) values that were corrected alongside the trial index update.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    for ($row = 2; $row -le 41; $row++) {
        $cell = $ws.Cells.Item($row, 2)
        $current = $cell.Value2
        $cell.Value2 = $current + 1
    }
}

# Sheet-specific prediction (column C) corrections
$wsP01 = $wb.Worksheets.Item("P01")
$wsP01.Cells.Item(9, 3).Value2 = 2

$wsP07 = $wb.Worksheets.Item("P07")
$wsP07.Cells.Item(2, 3).Value2 = 1

$wsP08 = $wb.Worksheets.Item("P08")
$wsP08.Cells.Item(26, 3).Value2 = 1
